$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("D2").Value = "56.435.16"
$ws.Range("E2").Value = "  -3.04%  "
$ws.Range("D3").Value = "2.990.58"
$ws.Range("E4").Value = "  +0.14%  "
$ws.Range("D5").Value = "'497.35"
$ws.Range("E5").Value = "  -5.22%  "
$ws.Range("D6").Value = "'135.32"
$ws.Range("E6").Value = "  +0.72%  "
$ws.Range("E7").Value = "  -0.03%  "
$ws.Range("E8").Value = "  -3.07%  "
$ws.Range("D9").Value = "'7.26"
$ws.Range("E9").Value = "  +1.29%  "
$ws.Range("E10").Value = "  -2.81%  "
$ws.Range("D11").Value = "'0.353"
$ws.Range("E11").Value = "  -6.57%  "
$ws.Range("E12").Value = "  -0.13%  "
$ws.Range("D13").Value = "3.502.49"
$ws.Range("E13").Value = "  -4.49%  "
$ws.Range("D14").Value = "'25.05"
$ws.Range("E14").Value = "  -1.59%  "
$ws.Range("D15").Value = "56.481.05"
$ws.Range("E15").Value = "  -2.87%  "
$ws.Range("D16").Value = "2.989.44"
$ws.Range("E16").Value = "  -4.77%  "
$ws.Range("E17").Value = "  -3.84%  "
$ws.Range("D18").Value = "'5.85"
$ws.Range("E18").Value = "  +1.49%  "
$ws.Range("D19").Value = "'12.40"
$ws.Range("E19").Value = "  -4.65%  "
$ws.Range("D20").Value = "'7.76"
$ws.Range("E20").Value = "  -1.74%  "
$ws.Range("D21").Value = "'325.70"
$ws.Range("E21").Value = "  -4.94%  "
$ws.Range("D22").Value = "'1.00"
$ws.Range("E22").Value = "  +0.18%  "
$ws.Range("E23").Value = "  -7.21%  "
$ws.Range("D24").Value = "'61.33"
$ws.Range("E24").Value = "  -9.35%  "
$ws.Range("E25").Value = "  +0.90%  "
$ws.Range("D26").Value = "'0.162"
$ws.Range("E26").Value = "  -3.18%  "
$ws.Range("D27").Value = "0.0₃0899"
$ws.Range("E27").Value = "  -6.02%  "
$ws.Range("E28").Value = "  -0.07%  "
$ws.Range("D29").Value = "'6.50"
$ws.Range("E29").Value = "  -3.89%  "
$ws.Range("D30").Value = "'6.67"
$ws.Range("E30").Value = "  -2.18%  "
$ws.Range("B31").Value = "PancakeSwap"
$ws.Range("C31").Value = "https://coinranking.com/coin/ncYFcP709+pancakeswap-cake"
$ws.Range("D31").Value = "'1.74"
$ws.Range("E31").Value = "  -6.21%  "
$ws.Range("B32").Value = "Fetch.AI"
$ws.Range("C32").Value = "https://coinranking.com/coin/AWma-WzFHmKVQ+fetchai-fet"
$ws.Range("D32").Value = "'1.17"
$ws.Range("E32").Value = "  -3.63%  "
$ws.Range("D33").Value = "'20.21"
$ws.Range("E33").Value = "  -5.44%  "
$ws.Range("D34").Value = "'155.47"
$ws.Range("E34").Value = "  -0.80%  "
$ws.Range("D35").Value = "'4.47"
$ws.Range("E35").Value = "  -6.31%  "
$ws.Range("D36").Value = "'1.28"
$ws.Range("E36").Value = "  -6.06%  "
$ws.Range("E37").Value = "  -9.58%  "
$ws.Range("D38").Value = "'0.0687"
$ws.Range("E38").Value = "  +0.43%  "
$ws.Range("D39").Value = "'23.28"
$ws.Range("E39").Value = "  -3.50%  "
$ws.Range("D40").Value = "3.024.20"
$ws.Range("E40").Value = "  -4.35%  "
$ws.Range("D41").Value = "'36.56"
$ws.Range("E41").Value = "  -9.34%  "
$ws.Range("E42").Value = "  -0.02%  "
$ws.Range("D43").Value = "'0.640"
$ws.Range("E43").Value = "  -7.21%  "
$ws.Range("B44").Value = "Maker"
$ws.Range("C44").Value = "https://coinranking.com/coin/qFakph2rpuMOL+maker-mkr"
$ws.Range("D44").Value = "2.231.36"
$ws.Range("E44").Value = "  -1.22%  "
$ws.Range("B45").Value = "ONDO"
$ws.Range("C45").Value = "https://coinranking.com/coin/7AQlxzQpQ+ondo-ondo"
$ws.Range("D45").Value = "'0.991"
$ws.Range("E45").Value = "  -7.92%  "
$ws.Range("B46").Value = "Stacks"
$ws.Range("C46").Value = "https://coinranking.com/coin/mMPrMcB7+stacks-stx"
$ws.Range("D46").Value = "'1.40"
$ws.Range("E46").Value = "  -2.16%  "
$ws.Range("D47").Value = "'3.57"
$ws.Range("E47").Value = "  -7.88%  "
$ws.Range("D48").Value = "'1.94"
$ws.Range("E48").Value = "  +5.06%  "
$ws.Range("E49").Value = "  -6.14%  "
$ws.Range("E50").Value = "  +1.18%  "
$ws.Range("D51").Value = "'19.06"
$ws.Range("E51").Value = "  -7.37%  "
